$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
# Clear the existing multi-run / line-break title text entirely, then set
# the new single-run title so no stray runs, <a:br/>, or endParaRPr remain.
$tf.DeleteText()
$tf.TextRange.Text = "Anomaly detection in graphs - past, present and future."
